$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like text value (e.g. "2023-09-05") into a cell
# without letting Excel auto-convert it to a date serial number.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# New records appended as rows 25-31 (columns follow the existing header
# layout: A Id, B Taxonsorteringsordning, C Valideringsstatus, D Rodlistade,
# E TaxonId, F Artnamn, G Vetenskapligt namn, H Auktor, P Lokalnamn, Q Ost,
# R Nord, S Noggrannhet, T Lan, U Kommun, V Provins, W Forsamling,
# Y Startdatum, Z Starttid, AA Slutdatum, AB Sluttid, AC Publik kommentar,
# AD Ej aterfunnen, AE Osaker artbestamning, AG Ospontan, AW Rapportor,
# AX Observatorer)

$rows = @(
    [ordered]@{
        Row=25; A=111936854; B=56414; C="Ovaliderad"; D="NT"; E=100049
        F="Spillkråka"; G="Dryocopus martius"; H="(Linnaeus, 1758)"
        P="rörvattsbodarna, Jmt"; Q=450998.3386916541; R=7087288.958247212; S=10
        T="Jämtland"; U="Krokom"; V="Jämtland"; W="Hotagen"
        Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00"
        AD=$false; AE=$false; AG=$false
        AW="Benny Öwre"; AX="Benny Öwre"
    }
    [ordered]@{
        Row=26; A=111936892; B=77515; C="Ovaliderad"; D="NT"; E=6425
        F="Garnlav"; G="Alectoria sarmentosa"; H="(Ach.) Ach."
        P="rörvattsbodarna, Jmt"; Q=451172.0902361136; R=7086726.569319103; S=10
        T="Jämtland"; U="Krokom"; V="Jämtland"; W="Hotagen"
        Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00"
        AD=$false; AE=$false; AG=$false
        AW="Benny Öwre"; AX="Benny Öwre"
    }
    [ordered]@{
        Row=27; A=111936864; B=89423; C="Ovaliderad"; D="NT"; E=5432
        F="Granticka"; G="Porodaedalea chrysoloma"; H="(Fr.) Fiasson & Niemelä"
        P="rörvattsbodarna, Jmt"; Q=451094.1385684713; R=7087212.607717684; S=10
        T="Jämtland"; U="Krokom"; V="Jämtland"; W="Hotagen"
        Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00"
        AD=$false; AE=$false; AG=$false
        AW="Benny Öwre"; AX="Benny Öwre"
    }
    [ordered]@{
        Row=28; A=111936894; B=89965; C="Ovaliderad"; D="VU"; E=760
        F="Doftticka"; G="Haploporus odorus"; H="(Sommerf.) Bondartsev & Singer"
        P="rörvattsbodarna, Jmt"; Q=451168.6101546783; R=7086616.526546557; S=10
        T="Jämtland"; U="Krokom"; V="Jämtland"; W="Hotagen"
        Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00"
        AD=$false; AE=$false; AG=$false
        AW="Benny Öwre"; AX="Benny Öwre"
    }
    [ordered]@{
        Row=29; A=111936789; B=90087; C="Ovaliderad"; D="LC"; E=3298
        F="Trådticka"; G="Climacocystis borealis"; H="(Fr.) Kotl. & Pouzar"
        P="rörvattsbodarna, Jmt"; Q=450955.1314140605; R=7087063.751596102; S=10
        T="Jämtland"; U="Krokom"; V="Jämtland"; W="Hotagen"
        Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00"
        AD=$false; AE=$false; AG=$false
        AW="Benny Öwre"; AX="Benny Öwre"
    }
    [ordered]@{
        Row=30; A=111936793; B=56398; C="Ovaliderad"; D="NT"; E=100109
        F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"
        P="rörvattsbodarna, Jmt"; Q=451088.7179698629; R=7087232.506422138; S=10
        T="Jämtland"; U="Krokom"; V="Jämtland"; W="Hotagen"
        Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00"
        AC="ringhack äldre"
        AD=$false; AE=$false; AG=$false
        AW="Benny Öwre"; AX="Benny Öwre"
    }
    [ordered]@{
        Row=31; A=111936860; B=89423; C="Ovaliderad"; D="NT"; E=5432
        F="Granticka"; G="Porodaedalea chrysoloma"; H="(Fr.) Fiasson & Niemelä"
        P="rörvattsbodarna, Jmt"; Q=450975.0281813644; R=7086982.857739178; S=10
        T="Jämtland"; U="Krokom"; V="Jämtland"; W="Hotagen"
        Y="2023-09-05"; Z="00:00"; AA="2023-09-05"; AB="00:00"
        AD=$false; AE=$false; AG=$false
        AW="Benny Öwre"; AX="Benny Öwre"
    }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H

    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
    $ws.Range("W$row").Value = $r.W

    Set-TextValue "Y$row" $r.Y
    $ws.Range("Z$row").Value = $r.Z
    Set-TextValue "AA$row" $r.AA
    $ws.Range("AB$row").Value = $r.AB

    if ($r.Contains("AC")) {
        $ws.Range("AC$row").Value = $r.AC
    }

    $ws.Range("AD$row").Value = $r.AD
    $ws.Range("AE$row").Value = $r.AE
    $ws.Range("AG$row").Value = $r.AG

    $ws.Range("AW$row").Value = $r.AW
    $ws.Range("AX$row").Value = $r.AX
}
